# Apply cryptos list update (prices/volumes refreshed; three coin rows re-ranked).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.750.88'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.295.31'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.78'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.31'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.09'
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.58'
$ws.Range('E12').Value = '  +4.85%  '
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.85'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '2.652.87'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '2.291.81'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.776'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = '42.671.24'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.69'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.02'
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.07'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.65'
$ws.Range('E23').Value = '  -2.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '166.48'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.05'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.01'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.78'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.91'
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('E34').Value = '  -0.61%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.36'
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  -8.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0686'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').Value = '1.994.70'
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.25'
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.17'
$ws.Range('E45').Value = '  +5.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.12'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.83'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.33'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.520.03'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.85'
$ws.Range('E51').Value = '  -1.71%  '
